$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 10:33"

# Row 7 - Rusia
$ws.Cells.Item(7, 2).Value = 1415316
$ws.Cells.Item(7, 3).Value = 15982
$ws.Cells.Item(7, 4).Value = 1075904
$ws.Cells.Item(7, 5).Value = 315046
$ws.Cells.Item(7, 7).Value = 179
$ws.Cells.Item(7, 8).Value = 24366

# Row 22 - Indonesia
$ws.Cells.Item(22, 2).Value = 365240
$ws.Cells.Item(22, 3).Value = 3373
$ws.Cells.Item(22, 4).Value = 289243
$ws.Cells.Item(22, 5).Value = 63380
$ws.Cells.Item(22, 7).Value = 106
$ws.Cells.Item(22, 8).Value = 12617

# Row 23 - Filipinas
$ws.Cells.Item(23, 2).Value = 359169
$ws.Cells.Item(23, 3).Value = 2638
$ws.Cells.Item(23, 4).Value = 310303
$ws.Cells.Item(23, 5).Value = 42191
$ws.Cells.Item(23, 7).Value = 26
$ws.Cells.Item(23, 8).Value = 6675

# Row 33 - Polonia
$ws.Cells.Item(33, 4).Value = 94014
$ws.Cells.Item(33, 5).Value = 78179

# Row 65 - Singapur
$ws.Cells.Item(65, 2).Value = 57915
$ws.Cells.Item(65, 3).Value = 4
$ws.Cells.Item(65, 5).Value = 80

# Row 84 - Eslovaquia
$ws.Cells.Item(84, 2).Value = 30695
$ws.Cells.Item(84, 3).Value = 860
$ws.Cells.Item(84, 4).Value = 7536
$ws.Cells.Item(84, 5).Value = 23067
$ws.Cells.Item(84, 7).Value = 4
$ws.Cells.Item(84, 8).Value = 92

# Row 87 - Croacia
$ws.Cells.Item(87, 2).Value = 25973
$ws.Cells.Item(87, 3).Value = 393
$ws.Cells.Item(87, 4).Value = 20529
$ws.Cells.Item(87, 5).Value = 5070
$ws.Cells.Item(87, 7).Value = 11
$ws.Cells.Item(87, 8).Value = 374

# Row 94 - Georgia
$ws.Cells.Item(94, 2).Value = 18663
$ws.Cells.Item(94, 3).Value = 1186
$ws.Cells.Item(94, 4).Value = 8338
$ws.Cells.Item(94, 5).Value = 10182
$ws.Cells.Item(94, 7).Value = 7
$ws.Cells.Item(94, 8).Value = 143

# Rows 118-120: Lituania moved above Angola/Mauritania (reorder + update)
# Row 118 becomes Lituania with new/updated figures
$ws.Cells.Item(118, 1).Value = "Lituania"
$ws.Cells.Item(118, 2).Value = 7726
$ws.Cells.Item(118, 3).Value = 205
$ws.Cells.Item(118, 4).Value = 3110
$ws.Cells.Item(118, 5).Value = 4503
$ws.Cells.Item(118, 8).Value = 113

# Row 119 becomes Angola (previous Angola figures, unchanged)
$ws.Cells.Item(119, 1).Value = "Angola"
$ws.Cells.Item(119, 2).Value = 7622
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 4).Value = 3030
$ws.Cells.Item(119, 5).Value = 4345
$ws.Cells.Item(119, 8).Value = 247

# Row 120 becomes Mauritania (previous Mauritania figures, unchanged)
$ws.Cells.Item(120, 1).Value = "Mauritania"
$ws.Cells.Item(120, 2).Value = 7608
$ws.Cells.Item(120, 3).Value = 0
$ws.Cells.Item(120, 4).Value = 7347
$ws.Cells.Item(120, 5).Value = 98
$ws.Cells.Item(120, 8).Value = 163

# Row 121 (Guadalupe) stays the same - no change needed

# Row 141 - Estonia
$ws.Cells.Item(141, 2).Value = 4085
$ws.Cells.Item(141, 3).Value = 7
$ws.Cells.Item(141, 4).Value = 3229
$ws.Cells.Item(141, 5).Value = 788

# Row 149 - Letonia
$ws.Cells.Item(149, 2).Value = 3494
$ws.Cells.Item(149, 3).Value = 44
$ws.Cells.Item(149, 4).Value = 1341
$ws.Cells.Item(149, 5).Value = 2109

$wb.Save()
